$d = $word.ActiveDocument

$d.Content.Find.Execute('111×5=555', $true, $false, $false, $false, $false, $true, 1, $false, '670×8=5360', 2) | Out-Null
$d.Content.Find.Execute('841×8=6728', $true, $false, $false, $false, $false, $true, 1, $false, '435×4=1740', 2) | Out-Null
$d.Content.Find.Execute('680×5=3400', $true, $false, $false, $false, $false, $true, 1, $false, '651×2=1302', 2) | Out-Null
$d.Content.Find.Execute('274×7=1918', $true, $false, $false, $false, $false, $true, 1, $false, '218×7=1526', 2) | Out-Null
$d.Content.Find.Execute('942×6=5652', $true, $false, $false, $false, $false, $true, 1, $false, '835×5=4175', 2) | Out-Null
$d.Content.Find.Execute('149×8=1192', $true, $false, $false, $false, $false, $true, 1, $false, '265×5=1325', 2) | Out-Null
$d.Content.Find.Execute('770×8=6160', $true, $false, $false, $false, $false, $true, 1, $false, '222×7=1554', 2) | Out-Null
$d.Content.Find.Execute('619×5=3095', $true, $false, $false, $false, $false, $true, 1, $false, '477×8=3816', 2) | Out-Null
$d.Content.Find.Execute('872×9=7848', $true, $false, $false, $false, $false, $true, 1, $false, '371×8=2968', 2) | Out-Null
$d.Content.Find.Execute('579×6=3474', $true, $false, $false, $false, $false, $true, 1, $false, '486×5=2430', 2) | Out-Null
$d.Content.Find.Execute('855×3=2565', $true, $false, $false, $false, $false, $true, 1, $false, '688×8=5504', 2) | Out-Null
$d.Content.Find.Execute('936×3=2808', $true, $false, $false, $false, $false, $true, 1, $false, '655×6=3930', 2) | Out-Null
$d.Content.Find.Execute('732×8=5856', $true, $false, $false, $false, $false, $true, 1, $false, '192×5=960', 2) | Out-Null
$d.Content.Find.Execute('696×6=4176', $true, $false, $false, $false, $false, $true, 1, $false, '454×6=2724', 2) | Out-Null
$d.Content.Find.Execute('612×8=4896', $true, $false, $false, $false, $false, $true, 1, $false, '542×7=3794', 2) | Out-Null
$d.Content.Find.Execute('969×2=1938', $true, $false, $false, $false, $false, $true, 1, $false, '786×6=4716', 2) | Out-Null
$d.Content.Find.Execute('662×3=1986', $true, $false, $false, $false, $false, $true, 1, $false, '244×4=976', 2) | Out-Null
$d.Content.Find.Execute('860×5=4300', $true, $false, $false, $false, $false, $true, 1, $false, '892×8=7136', 2) | Out-Null
$d.Content.Find.Execute('698×6=4188', $true, $false, $false, $false, $false, $true, 1, $false, '410×4=1640', 2) | Out-Null
$d.Content.Find.Execute('672×9=6048', $true, $false, $false, $false, $false, $true, 1, $false, '446×2=892', 2) | Out-Null
$d.Content.Find.Execute('925×5=4625', $true, $false, $false, $false, $false, $true, 1, $false, '277×2=554', 2) | Out-Null
$d.Content.Find.Execute('178×2=356', $true, $false, $false, $false, $false, $true, 1, $false, '236×8=1888', 2) | Out-Null
$d.Content.Find.Execute('335×8=2680', $true, $false, $false, $false, $false, $true, 1, $false, '307×2=614', 2) | Out-Null
$d.Content.Find.Execute('721×5=3605', $true, $false, $false, $false, $false, $true, 1, $false, '335×3=1005', 2) | Out-Null
$d.Content.Find.Execute('155×4=620', $true, $false, $false, $false, $false, $true, 1, $false, '292×7=2044', 2) | Out-Null
